$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Normalize the header row's box border first: with a 4th column being
# added, the old thick outer box (mixed thin/medium edges per cell) becomes
# a plain thin border around every header cell instead.
# ---------------------------------------------------------------------------
$ws.Range("C1").Borders.LineStyle = 1
$ws.Range("C1").Borders.Weight = 2
$ws.Range("A1").Borders(9).LineStyle = -4142

# Give B1 the same (now-thin-boxed) bold/shaded header look as C1.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Add the new "common_password file" column (D): the wordlist that would
# crack each row's password.
# ---------------------------------------------------------------------------
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "common_password file"

$ws.Range("A3").Copy()
$ws.Range("D2:D7").PasteSpecial(-4122)
$ws.Range("D2").Value = "common_passwords1.txt"
$ws.Range("D3").Value = "common_passwords2.txt"
$ws.Range("D4").Value = "common_passwords3.txt"
$ws.Range("D5").Value = "common_passwords1.txt"
$ws.Range("D6").Value = "common_passwords3.txt"
$ws.Range("D7").Value = "common_passwords1.txt"

# ---------------------------------------------------------------------------
# Row 2 (first data row) drops its heavier top/fill styling and matches the
# plain look used by the rest of the data rows.
# ---------------------------------------------------------------------------
$ws.Range("A3").Copy()
$ws.Range("A2:B2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Restore the default (auto) row height / drop the thick-bottom flag on the
# header row, and size the newly-populated columns C & D to fit their text,
# same as the existing bestFit columns A and B.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()

# ---------------------------------------------------------------------------
# Selection, as left by the author after the edit.
# ---------------------------------------------------------------------------
$ws.Range("D10").Select()
